# build based on 0294d82
# Re-run of the model without the "generator" asset: drops x_us_generator /
# CAPEX_us_generator / C_OEM_us_generator columns from design_users and
# economics_users, and refreshes the recomputed numeric results across all
# sheets.

$wb = $excel.ActiveWorkbook

$wsInfo       = $wb.Worksheets.Item("info_solution")
$wsEconAgg    = $wb.Worksheets.Item("economics_aggregator")
$wsPeakAgg    = $wb.Worksheets.Item("peak_aggregator")
$wsDesign     = $wb.Worksheets.Item("design_users")
$wsEconUsers  = $wb.Worksheets.Item("economics_users")
$wsPeakUsers  = $wb.Worksheets.Item("peak_users")

# --- structural change: remove the "generator" user-asset columns ---
# design_users: column F is x_us_generator
$wsDesign.Columns("F").Delete()

# economics_users: column M is CAPEX_us_generator, column S is
# C_OEM_us_generator. Delete the higher-indexed column first so the
# second delete still targets the right (un-shifted) column letter.
$wsEconUsers.Columns("S").Delete()
$wsEconUsers.Columns("M").Delete()

# --- refreshed numeric results from the re-run ---
$wsInfo.Range("A2").Value = 1.3108773231506348

$wsEconAgg.Range("A2").Value = 8530.582818480752
$wsEconAgg.Range("C2").Value = 8530.582818480752
$wsEconAgg.Range("D2").Value = -1214211.3113248872
$wsEconAgg.Range("E2").Value = -1222741.894143368
$wsEconAgg.Range("F2").Value = 573.3891603577929

$wsPeakAgg.Range("B2").Value = 74.95619728130688
$wsPeakAgg.Range("C2").Value = 103.67293866601534
$wsPeakAgg.Range("F2").Value = 63.832102472583365
$wsPeakAgg.Range("G2").Value = 71.09777497045125
$wsPeakAgg.Range("H2").Value = 83.61572241881349
$wsPeakAgg.Range("I2").Value = 75.55808313800486
$wsPeakAgg.Range("J2").Value = 71.05252151118073
$wsPeakAgg.Range("K2").Value = 67.73624642560675
$wsPeakAgg.Range("L2").Value = 78.81711346915316

$wsDesign.Range("E4").Value = 44.04652629600474
$wsDesign.Range("F4").Value = 3.7213684553918602
$wsDesign.Range("G4").Value = 3.7213684553918602
$wsDesign.Range("H4").Value = 34.906983718279804

$wsEconUsers.Range("E2").Value = 0.0
$wsEconUsers.Range("B3").Value = -223241.12057109867
$wsEconUsers.Range("D3").Value = -12164.328747674772
$wsEconUsers.Range("J3").Value = -169392.5550729593
$wsEconUsers.Range("C4").Value = 177428.2143016821
$wsEconUsers.Range("D4").Value = -17171.222950664825
$wsEconUsers.Range("F4").Value = 35626.4177709583
$wsEconUsers.Range("G4").Value = 1509.251446052555
$wsEconUsers.Range("I4").Value = 9783.233210352735
$wsEconUsers.Range("J4").Value = -219838.02000083425
$wsEconUsers.Range("L4").Value = 70474.44207360758
$wsEconUsers.Range("M4").Value = 1488.547382156744
$wsEconUsers.Range("N4").Value = 744.273691078372
$wsEconUsers.Range("O4").Value = 104720.95115483941
$wsEconUsers.Range("Q4").Value = 1321.3957888801424
$wsEconUsers.Range("R4").Value = 18.606842276959302
$wsEconUsers.Range("S4").Value = 7.4427369107837205
$wsEconUsers.Range("T4").Value = 1047.2095115483942

$wsPeakUsers.Range("O2").Value = 32.86113926823138
